# Applies the edits described by the commit diff to the "Data" worksheet
# (sheet2.xml in the package): updated T-column kinetic values for a few
# existing rows, newly populated rows 33-42 (previously blank template
# rows), and an updated selection/active-cell on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

# --- Updated values in column T (existing rows) ---
$ws.Range("T3").Value = 0.4
$ws.Range("T4").Value = 0.5
$ws.Range("T5").Value = 0.5
$ws.Range("T6").Value = 0.9
$ws.Range("T24").Value = 0.9

# --- Newly populated data rows 33-42 ---
# Column A: temperature-ish series, column H mirrors A.
$aValues = @(20, 25, 30, 40, 50, 60, 70, 80, 90, 100)
$tValues = @(1.5, 5.5, 9.4, 11, 12, 12.3, 12.8, 12.7, 12.9, 11.8)

for ($i = 0; $i -lt 10; $i++) {
    $row = 33 + $i
    $ws.Range("A$row").Value = $aValues[$i]
    $ws.Range("B$row").Value = 160
    $ws.Range("C$row").Value = 9
    $ws.Range("E$row").Value = 0.5
    $ws.Range("F$row").Value = 15
    $ws.Range("G$row").Value = 12
    $ws.Range("H$row").Value = $aValues[$i]
    $ws.Range("I$row").Value = 0
    $ws.Range("J$row").Value = 0
    $ws.Range("T$row").Value = $tValues[$i]
    $ws.Range("Y$row").Value = 0
}

# Column D carries the shared "moles of acid" formula down through row 42.
$ws.Range("D33:D42").Formula = "=2000*0.65/100/98.079"

# Column N carries the shared dilution-factor formula down through row 42.
$ws.Range("N33:N42").Formula = "=15.67/0.88"

# --- Updated window scroll position / selection on the Data sheet ---
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 9
$ws.Range("V29").Select()
